$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reg")

$ws.Range("C2").Value = "akhilbingi321312"
$ws.Range("C3").Value = "akhilbingi642344423"
$ws.Range("C4").Value = "akhilbingi32445534"
$ws.Range("C5").Value = "akhilbingi32156645"
$ws.Range("C6").Value = "akhilbingi321277756"
$ws.Range("C7").Value = "nffmf8867"
$ws.Range("C8").Value = "akhilhdiw59978"
$ws.Range("C9").Value = "akhilhdingi90089"
$ws.Range("C10").Value = "RoyalEnfield1"

$ws.Range("C10").Select()
